$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 165, shifting existing rows 165-181 down to 166-182
$ws.Rows.Item(165).Insert()

# Populate the newly inserted row 165 with the new data point
$ws.Cells.Item(165, 1).Value = 11
$ws.Cells.Item(165, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(165, 3).Value = "Bíobío"
$ws.Cells.Item(165, 4).Value = Get-Date -Year 2023 -Month 7 -Day 25 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(165, 5).Value = 8
$ws.Cells.Item(165, 6).Value = "Fruta"
$ws.Cells.Item(165, 7).Value = 100108
$ws.Cells.Item(165, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(165, 9).Value = 100108002
$ws.Cells.Item(165, 10).Value = "Mango"
$ws.Cells.Item(165, 11).Value = "Sin especificar"
$ws.Cells.Item(165, 12).Value = "Primera"
$ws.Cells.Item(165, 13).Value = 150
$ws.Cells.Item(165, 14).Value = 8000
$ws.Cells.Item(165, 15).Value = 8000
$ws.Cells.Item(165, 16).Value = 8000
$ws.Cells.Item(165, 17).Value = '$/bandeja 4 kilos'
$ws.Cells.Item(165, 18).Value = "Perú"
$ws.Cells.Item(165, 19).Value = 2000
$ws.Cells.Item(165, 20).Value = 4
